$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewicht")

# --- 1. Make room for a new row 2 (source citation) + keep row 3 blank ---
# (insert twice, clearing the copied-down format each time so the new
#  rows come out as plain/unstyled, matching a never-formatted row)
$ws.Range("A2").EntireRow.Insert()
$ws.Range("A2").ClearFormats()
$ws.Range("A2").EntireRow.Insert()
$ws.Range("A2").ClearFormats()

# --- 2. Insert 4 new rows for the "3-17" age groups, Male block ---
# (old row 4 "18 und mehr"/Male is now row 6)
$ws.Range("A6:A9").EntireRow.Insert()

# --- 3. Insert 4 new rows for the "3-17" age groups, Female block ---
# (old row 12 "18 und mehr"/Female is now row 18 after the prior inserts)
$ws.Range("A18:A21").EntireRow.Insert()

# --- 4. Fill the new Male rows (6-9) ---
$ws.Range("A6").Value = "3 - 6"
$ws.Range("B6").Value = "Male"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 86.7
$ws.Range("E6").Value = 6.4
$ws.Range("F6").Value = 1
$ws.Range("G6").Formula = "=SUM(C6:F6)"

$ws.Range("A7").Value = "7 - 10"
$ws.Range("B7").Value = "Male"
$ws.Range("C7").Value = 7.6
$ws.Range("D7").Value = 76.3
$ws.Range("E7").Value = 9.3
$ws.Range("F7").Value = 6.8
$ws.Range("G7").Formula = "=SUM(C7:F7)"

$ws.Range("A8").Value = "11 - 13"
$ws.Range("B8").Value = "Male"
$ws.Range("C8").Value = 8.5
$ws.Range("D8").Value = 70.4
$ws.Range("E8").Value = 13.1
$ws.Range("F8").Value = 8
$ws.Range("G8").Formula = "=SUM(C8:F8)"

$ws.Range("A9").Value = "14 - 17"
$ws.Range("B9").Value = "Male"
$ws.Range("C9").Value = 9.7
$ws.Range("D9").Value = 71.7
$ws.Range("E9").Value = 9.3
$ws.Range("F9").Value = 9.2
$ws.Range("G9").Formula = "=SUM(C9:F9)"

# --- 5. Fill the new Female rows (18-21) ---
$ws.Range("A18").Value = "3 - 6"
$ws.Range("B18").Value = "Female"
$ws.Range("C18").Value = 4.3
$ws.Range("D18").Value = 84.8
$ws.Range("E18").Value = 7.6
$ws.Range("F18").Value = 3.2
$ws.Range("G18").Formula = "=SUM(C18:F18)"

$ws.Range("A19").Value = "7 - 10"
$ws.Range("B19").Value = "Female"
$ws.Range("C19").Value = 9.6
$ws.Range("D19").Value = 75.5
$ws.Range("E19").Value = 10.3
$ws.Range("F19").Value = 4.7
$ws.Range("G19").Formula = "=SUM(C19:F19)"

$ws.Range("A20").Value = "11 - 13"
$ws.Range("B20").Value = "Female"
$ws.Range("C20").Value = 7.7
$ws.Range("D20").Value = 72.3
$ws.Range("E20").Value = 13.5
$ws.Range("F20").Value = 6.5
$ws.Range("G20").Formula = "=SUM(C20:F20)"

$ws.Range("A21").Value = "14 - 17"
$ws.Range("B21").Value = "Female"
$ws.Range("C21").Value = 7.3
$ws.Range("D21").Value = 76.5
$ws.Range("E21").Value = 8.5
$ws.Range("F21").Value = 7.7
$ws.Range("G21").Formula = "=SUM(C21:F21)"

# --- 6. Mark the new age-group labels as Text so they don't get
#        reinterpreted (e.g. "3 - 6" looking like a date/range) ---
$ws.Range("A6:A9").NumberFormat = "@"
$ws.Range("A18:A21").NumberFormat = "@"

# --- 7. Add the data-source citation in row 2 (added last so the new
#        shared strings land in the same order as upstream) ---
$ws.Range("A2").Value = "Tab 1 in https://www.rki.de/DE/Content/Kommissionen/Bundesgesundheitsblatt/Downloads/2019_10_Schienkiewitz_BMI.pdf?__blob=publicationFile"

# --- 8. Make "Gewicht" the active/selected sheet + cell, like the author left it ---
$ws.Activate()
$ws.Range("C6").Select() | Out-Null
